$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row 2
Set-TextCell 2 4 "71.382.37"
Set-TextCell 2 5 "  +0.12%  "

# Row 3
Set-TextCell 3 4 "3.821.43"
Set-TextCell 3 5 "  -0.50%  "

# Row 4
Set-TextCell 4 4 "0.999"
Set-TextCell 4 5 "  -0.02%  "

# Row 5
Set-TextCell 5 4 "704.10"
Set-TextCell 5 5 "  -1.58%  "

# Row 6
Set-TextCell 6 4 "171.31"
Set-TextCell 6 5 "  -1.03%  "

# Row 7
Set-TextCell 7 4 "3.820.49"
Set-TextCell 7 5 "  -0.40%  "

# Row 8
Set-TextCell 8 5 "  -0.08%  "

# Row 9
Set-TextCell 9 5 "  -0.21%  "

# Row 10
Set-TextCell 10 5 "  -1.90%  "

# Row 11
Set-TextCell 11 4 "7.50"
Set-TextCell 11 5 "  +1.89%  "

# Row 12
Set-TextCell 12 4 "0.488"
Set-TextCell 12 5 "  +5.94%  "

# Row 13
Set-TextCell 13 4 "0.0000252"
Set-TextCell 13 5 "  -1.92%  "

# Row 14
Set-TextCell 14 4 "36.78"
Set-TextCell 14 5 "  -0.21%  "

# Row 15
Set-TextCell 15 4 "4.461.68"
Set-TextCell 15 5 "  -0.64%  "

# Row 16
Set-TextCell 16 4 "3.807.28"
Set-TextCell 16 5 "  -0.80%  "

# Row 17
Set-TextCell 17 4 "71.501.30"
Set-TextCell 17 5 "  +0.35%  "

# Row 18
Set-TextCell 18 4 "7.25"
Set-TextCell 18 5 "  +0.08%  "

# Row 19
Set-TextCell 19 4 "17.60"
Set-TextCell 19 5 "  +0.93%  "

# Row 20
Set-TextCell 20 5 "  +0.18%  "

# Row 21
Set-TextCell 21 4 "513.52"
Set-TextCell 21 5 "  +3.61%  "

# Row 22
Set-TextCell 22 4 "10.54"
Set-TextCell 22 5 "  -1.73%  "

# Row 23
Set-TextCell 23 4 "0.718"
Set-TextCell 23 5 "  -1.53%  "

# Row 24
Set-TextCell 24 4 "83.96"
Set-TextCell 24 5 "  -1.57%  "

# Row 25
Set-TextCell 25 4 "0.0000142"
Set-TextCell 25 5 "  -2.85%  "

# Row 26
Set-TextCell 26 4 "12.79"
Set-TextCell 26 5 "  +4.91%  "

# Row 27
Set-TextCell 27 4 "3.961.75"
Set-TextCell 27 5 "  -0.79%  "

# Row 28
Set-TextCell 28 4 "10.37"
Set-TextCell 28 5 "  -3.04%  "

# Row 29
Set-TextCell 29 5 "  +0.04%  "

# Row 30
Set-TextCell 30 4 "2.00"

# Row 31
Set-TextCell 31 5 "  -5.23%  "

# Row 32
Set-TextCell 32 5 "  +1.24%  "

# Row 33
Set-TextCell 33 4 "7.39"
Set-TextCell 33 5 "  -1.39%  "

# Row 34
Set-TextCell 34 4 "29.31"
Set-TextCell 34 5 "  -0.45%  "

# Row 35
Set-TextCell 35 5 "  -4.97%  "

# Row 36
Set-TextCell 36 4 "9.35"
Set-TextCell 36 5 "  +1.38%  "

# Row 37
Set-TextCell 37 4 "3.782.34"
Set-TextCell 37 5 "  -0.62%  "

# Row 38
Set-TextCell 38 4 "0.998"
Set-TextCell 38 5 "  +0.21%  "

# Row 39
Set-TextCell 39 4 "6.69"
Set-TextCell 39 5 "  +11.02%  "

# Row 40
Set-TextCell 40 5 "  -1.71%  "

# Row 41
Set-TextCell 41 4 "2.39"
Set-TextCell 41 5 "  +4.69%  "

# Row 42
Set-TextCell 42 5 "  -2.11%  "

# Row 43
Set-TextCell 43 2 "dogwifhat"
Set-TextCell 43 3 "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextCell 43 4 "3.24"
Set-TextCell 43 5 "  -3.40%  "

# Row 44
Set-TextCell 44 2 "USDe"
Set-TextCell 44 3 "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextCell 44 4 "1.00"
Set-TextCell 44 5 "  +0.00%  "

# Row 45
Set-TextCell 45 5 "  +0.04%  "

# Row 46
Set-TextCell 46 4 "166.27"
Set-TextCell 46 5 "  +1.50%  "

# Row 47
Set-TextCell 47 4 "49.98"
Set-TextCell 47 5 "  +2.16%  "

# Row 48
Set-TextCell 48 4 "432.74"
Set-TextCell 48 5 "  +1.63%  "

# Row 49
Set-TextCell 49 5 "  -5.65%  "

# Row 50
Set-TextCell 50 4 "30.76"
Set-TextCell 50 5 "  +8.60%  "

# Row 51
Set-TextCell 51 2 "Cosmos"
Set-TextCell 51 3 "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextCell 51 4 "8.68"
Set-TextCell 51 5 "  -0.66%  "
